# Update "want to go" / 想去人数 counts (column F) on several rows across
# three worksheets: 展览 (Exhibitions), 本地生活 (Local life) and 全部类型 (All types).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 2884
$ws1.Range("F19").Value = 2169
$ws1.Range("F23").Value = 1056
$ws1.Range("F24").Value = 72
$ws1.Range("F28").Value = 778
$ws1.Range("F31").Value = 23
$ws1.Range("F37").Value = 367
$ws1.Range("F38").Value = 2386
$ws1.Range("F45").Value = 299
$ws1.Range("F46").Value = 111

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value  = 1840
$ws3.Range("F11").Value = 846

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 2884
$ws4.Range("F13").Value = 846
$ws4.Range("F20").Value = 2169
$ws4.Range("F25").Value = 1056
$ws4.Range("F26").Value = 72
$ws4.Range("F31").Value = 23
$ws4.Range("F38").Value = 367
$ws4.Range("F41").Value = 2386
$ws4.Range("F44").Value = 299
$ws4.Range("F45").Value = 111
